$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "lam xong danh muc ho so" - finish the record/profile catalog:
#   - add more dummy "tài khoản N" test accounts (rows 3, 10, 11 are new;
#     rows 4-9 and 12 get renumbered account names)
#   - reset the TrangThai (status, column E) flag back to 0 for every record
#   - move the last cell selection to H8
# ---------------------------------------------------------------------------

# Column A: account name labels for rows 2-12 ("tài khoản 1" .. "tài khoản 11")
$ws.Range("A2").Value  = "tài khoản 1"
$ws.Range("A3").Value  = "tài khoản 2"
$ws.Range("A4").Value  = "tài khoản 3"
$ws.Range("A5").Value  = "tài khoản 4"
$ws.Range("A6").Value  = "tài khoản 5"
$ws.Range("A7").Value  = "tài khoản 6"
$ws.Range("A8").Value  = "tài khoản 7"
$ws.Range("A9").Value  = "tài khoản 8"
$ws.Range("A10").Value = "tài khoản 9"
$ws.Range("A11").Value = "tài khoản 10"
$ws.Range("A12").Value = "tài khoản 11"

# Column E: TrangThai status reset to 0 for every record (rows 2-12)
$ws.Range("E2").Value  = 0
$ws.Range("E3").Value  = 0
$ws.Range("E4").Value  = 0
$ws.Range("E5").Value  = 0
$ws.Range("E6").Value  = 0
$ws.Range("E7").Value  = 0
$ws.Range("E8").Value  = 0
$ws.Range("E9").Value  = 0
$ws.Range("E10").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("E12").Value = 0

# Last UI selection moves to H8
$ws.Range("H8").Select()

Write-Host "done"
